# Fix mojibake: "Â±" (UTF-8 bytes for ± mis-decoded as Latin-1) -> "±"
# Affects columns B, C, D for rows 2 through 17 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badChar = [string][char]0x00C2 + [string][char]0x00B1   # "Â±"
$goodChar = [string][char]0x00B1                          # "±"

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $value = $cell.Value2
        if ($value -ne $null -and $value.ToString().Contains($badChar)) {
            $cell.Value = $value.ToString().Replace($badChar, $goodChar)
        }
    }
}
